# TFS 5404 - Allow users with job code starting "WPOP" to have access to
# Historical Dashboard. (formerly TFS 3878 - Email CSR comments text)

$d = $word.ActiveDocument

$dash = [char]0x2013
$lq   = [char]0x201c
$rq   = [char]0x201d

# ---------------------------------------------------------------------
# 1) Title-block "Description" cell (Table 1): replace the old TFS 3878
#    sentence with the new bold TFS 5404 sentence.
# ---------------------------------------------------------------------
$descCell = $d.Tables.Item(1).Rows.Item(1).Cells.Item(2)
$oldDesc = "3878 $dash Email CSR comments to supervisor/manager when CSR logs are completed."
$newDesc = "5404"

# NOTE: Replace must be wdReplaceOne (1), scoped to the cell's Range, so
# that only this single occurrence is touched (wdReplaceAll searches the
# whole story regardless of the Range the Find object was created on).
$f = $descCell.Range.Find
$f.ClearFormatting()
$f.Replacement.ClearFormatting()
$f.Execute($oldDesc, $false, $false, $false, $false, $false, $true, 1, $false, $newDesc, 1) | Out-Null

# Append the remaining bold text (" - Allow users with job codes starting
# "WPOP" to access Historical Dashboard.") right after "5404". Re-fetch the
# cell fresh from the document after the mutation above -- a Range handle
# captured before a Find/Replace edit keeps stale Start/End offsets.
$descCell = $d.Tables.Item(1).Rows.Item(1).Cells.Item(2)
$cellEnd = $descCell.Range.End
$insPt = $d.Range($cellEnd - 1, $cellEnd - 1)
$insPt.InsertAfter(" $dash Allow users with job codes starting " + $lq + "WPOP" + $rq + " to access Historical Dashboard.")
$insPt.Font.Bold = 1

# ---------------------------------------------------------------------
# 2) Revision-history table (Table 2): add a new last row documenting
#    the TFS 5404 change.
# ---------------------------------------------------------------------
$histTable = $d.Tables.Item(2)
$newRow = $histTable.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "2/1/2017"
$newRow.Cells.Item(2).Range.Text = "TFS 5404 $dash Allow users with job codes starting " + $lq + "WPOP" + $rq + " to access Historical Dashboard"
$newRow.Cells.Item(3).Range.Text = "Lili Huang"

# ---------------------------------------------------------------------
# 3) Narrative references to "TFS 3878." -> "TFS 5404."
# ---------------------------------------------------------------------
$r1 = $d.Content.Duplicate
$f1 = $r1.Find
$f1.ClearFormatting()
$f1.Replacement.ClearFormatting()
$f1.Execute("per TFS 3878.", $false, $false, $false, $false, $false, $true, 1, $false, "per TFS 5404.", 1) | Out-Null

$r2 = $d.Content.Duplicate
$f2 = $r2.Find
$f2.ClearFormatting()
$f2.Replacement.ClearFormatting()
$f2.Execute("for TFS 3878.", $false, $false, $false, $false, $false, $true, 1, $false, "for TFS 5404.", 1) | Out-Null

# ---------------------------------------------------------------------
# 4) Build artifact path: drop "_publish" from the zip file name and
#    bump the changeset number.
# ---------------------------------------------------------------------
$r3 = $d.Content.Duplicate
$f3 = $r3.Find
$f3.ClearFormatting()
$f3.Replacement.ClearFormatting()
$f3.Execute("eCoaching_publish.zip", $false, $false, $false, $false, $false, $true, 1, $false, "eCoaching.zip", 1) | Out-Null

$r4 = $d.Content.Duplicate
$f4 = $r4.Find
$f4.ClearFormatting()
$f4.Replacement.ClearFormatting()
$f4.Execute("C36194", $false, $false, $false, $false, $false, $true, 1, $false, "C36498", 1) | Out-Null

Write-Output "done"
